$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" – refresh the "Latest Handback DateTime"
# column for the just-handed-back file (row 2, the 2df024e5-... file) on
# both locale sheets.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-21 01:12:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-21 01:12:50"
